$d = $word.ActiveDocument

$result = $d.Content.Find.Execute("Unveiling the Complexity of Cosmic Phenomena", $true, $false, $false, $false, $false, $true, 1, $false, "Exploring the Marvelous World of Chemistry: A Journey into the Realm of Elements and Compounds", 2)
if (-not $result) { Write-Host "FAILED replace #0: Unveiling the Complexity of Cosmic Pheno" }
$result = $d.Content.Find.Execute(" Neil deGrasse Tyson", $true, $false, $false, $false, $false, $true, 1, $false, " Eleanor Stanton", 2)
if (-not $result) { Write-Host "FAILED replace #1:  Neil deGrasse Tyson" }
$result = $d.Content.Find.Execute("NeilTyson@SpaceInstitute", $true, $false, $false, $false, $false, $true, 1, $false, "estanton@edu", 2)
if (-not $result) { Write-Host "FAILED replace #2: NeilTyson@SpaceInstitute" }
$result = $d.Content.Find.Execute("Embarking on an enchanting voyage to unravel the complexities of cosmic phenomena, we delve into the profound mysteries that enchant our universe", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry, an intriguing and impactful science, unveils the hidden intricacies of matter and its diverse interactions", 2)
if (-not $result) { Write-Host "FAILED replace #3: Embarking on an enchanting voyage to unr" }
$result = $d.Content.Find.Execute(" From the grand tapestry of galaxies that stretch across unfathomable distances, to the enigmatic fabric of space and time warping around celestial bodies, our quest for comprehension leads us down a path of captivating discoveries", $true, $false, $false, $false, $false, $true, 1, $false, " From the vast universe to the microscopic realm within our bodies, chemistry plays a pivotal role in shaping our world", 2)
if (-not $result) { Write-Host "FAILED replace #4:  From the grand tapestry of galaxies tha" }
$result = $d.Content.Find.Execute(" Through meticulous observation, tireless calculations, and imaginative leap, humanity continues to unlock the secrets of the cosmos, redefining our perception of existence itself", $true, $false, $false, $false, $false, $true, 1, $false, " As we delve into the fascinating tapestry of chemistry, we embark on an exhilarating exploration of the elements that constitute everything around us and the myriad compounds formed through their intricate combinations", 2)
if (-not $result) { Write-Host "FAILED replace #5:  Through meticulous observation, tireles" }
$result = $d.Content.Find.Execute("As our telescopes peer deeper into the vast expanse, we encounter distant worlds that ignite our curiosity and contemplation", $true, $false, $false, $false, $false, $true, 1, $false, "In this realm of substances, we uncover the fundamental principles governing chemical reactions, witnessing the wondrous transformations of matter into new entities with unique properties", 2)
if (-not $result) { Write-Host "FAILED replace #6: As our telescopes peer deeper into the v" }
$result = $d.Content.Find.Execute(" The interplay of cosmic forces, the birth and death of stars, and the symphony of interactions between celestial bodies fuel our inquiry into the fundamental principles that govern the universe's evolution", $true, $false, $false, $false, $false, $true, 1, $false, " The symphony of chemistry encompasses myriad concepts, from atomic structures and bonding arrangements to energy transfer and reaction dynamics", 2)
if (-not $result) { Write-Host "FAILED replace #7:  The interplay of cosmic forces, the bir" }
$result = $d.Content.Find.Execute(" Each cosmic event, each intricate celestial dance, holds clues to unraveling mysteries that have captivated humankind for millennia, painting a breathtaking tapestry of cosmic beauty and awe", $true, $false, $false, $false, $false, $true, 1, $false, " Each element, with its distinctive characteristics, contributes to the intricate dance of chemical interactions, orchestrating the formation of countless compounds with diverse applications in fields ranging from medicine to materials science", 2)
if (-not $result) { Write-Host "FAILED replace #8:  Each cosmic event, each intricate celes" }
$result = $d.Content.Find.Execute("Yet, the complexities of the universe extend beyond the reaches of our tangible world, delving into realms that transcend our current understanding", $true, $false, $false, $false, $false, $true, 1, $false, "As we unravel the enigmas of chemistry, we gain invaluable insights into the natural world, unveiling the intricate mechanisms underlying life itself", 2)
if (-not $result) { Write-Host "FAILED replace #9: Yet, the complexities of the universe ex" }
$result = $d.Content.Find.Execute(" From the mysteries of dark matter and energy, whose enigmatic nature eludes our grasp, to the theoretical concept of multiple universes and the complexities of multi-dimensional space, our exploration into the cosmos opens doors to realms that challenge our conceptual boundaries", $true, $false, $false, $false, $false, $true, 1, $false, " From the intricate workings of photosynthesis, the process by which plants convert sunlight into energy, to the intricate pathways of cellular respiration, the fundamental energy-generating process within living organisms, chemistry unveils the symphony of life at its most fundamental level", 2)
if (-not $result) { Write-Host "FAILED replace #10:  From the mysteries of dark matter and e" }
$result = $d.Content.Find.Execute(" In this journey of seeking cosmic knowledge, we embark on an intellectual adventure where wonder and enigma intersect, inviting us to the depths of the universe's profound secrets", $true, $false, $false, $false, $false, $true, 1, $false, " Its principles permeate every aspect of our existence, shaping the materials we use, the medicines that heal us, and the intricate complexity of the living world", 2)
if (-not $result) { Write-Host "FAILED replace #11:  In this journey of seeking cosmic knowl" }
$result = $d.Content.Find.Execute("Our journey of exploration into cosmic phenomena unveils a harmonious tapestry of elegance, mystery, and boundless beauty", $true, $false, $false, $false, $false, $true, 1, $false, "This essay embarks on an enthralling exploration of chemistry, venturing into the captivating realm of elements, compounds, and their captivating interactions", 2)
if (-not $result) { Write-Host "FAILED replace #12: Our journey of exploration into cosmic p" }
$result = $d.Content.Find.Execute(" From the observable grandeur of distant galaxies and the intricacies of gravitational dance to the enigma of unseen forces and the complexities of spacetime, the universe constantly presents us with profound mysteries that test the limits of our comprehension", $true, $false, $false, $false, $false, $true, 1, $false, " From the fundamental principles governing chemical reactions to the intricacies of life itself, chemistry weaves the tapestry of our world, influencing countless aspects of our existence", 2)
if (-not $result) { Write-Host "FAILED replace #13:  From the observable grandeur of distant" }
$result = $d.Content.Find.Execute(" As we continue to probe the cosmos with unrelenting curiosity, we are reminded that the path to cosmic understanding is a never-ending adventure, filled with captivating discoveries and awe-inspiring revelations that paint a mesmerizing portrait of the intricate ", $true, $false, $false, $false, $false, $true, 1, $false, " Through the study of chemistry, we gain a profound understanding of the natural world, unlocking the secrets of matter and its remarkable transformations, revealing the symphony of life at its most fundamental level", 2)
if (-not $result) { Write-Host "FAILED replace #14:  As we continue to probe the cosmos with" }
$toDelete = "universe we inhabit. Embracing the enigma of the cosmos, we venture forth with an insatiable hunger for knowledge, eager to unravel the secrets that the universe holds"
$result = $d.Content.Find.Execute($toDelete, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $result) { Write-Host "FAILED deletion" }

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Name = "Times New Roman"
}

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
